$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colors_table")

# The "color_03.png" icon file was re-uploaded to the GitHub repo under a new
# commit hash and with the filename now containing a URL-encoded space
# ("color%2003.png") instead of an underscore ("color_03.png"). Update the
# Icon_url cell for row 3 (Magenta-Purple, E4) to point at the new location.
$ws.Range("E4").Value = "https://github.com/Ing-Aladar-Dukay/CV_Dukay/blob/46ceb0c7b57c0fa1196c0a1bf7357bb73838ec6b/03%20Colors%20icons/color%2003.png"
